# Add two new records (cw0204180, cw0204181) to Sheet1, plus supporting
# new columns (subject3, publisher, geographic_subject).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column before the existing "form" column (col F) -------
# This shifts form/date/date_edtf one column to the right (F->G, G->H, H->I)
# while preserving their existing per-column styles, and leaves the new
# column F empty/ready for the "subject3" header + data.
$ws.Columns.Item(6).Insert()

# --- Fill in the two new data rows (17 & 18) + the three new headers -----
# Values are written in this specific order so that new shared strings are
# appended to the shared-string table in the same sequence as the target.
$ws.Range("A17").Value = "cw0204180"
$ws.Range("A18").Value = "cw0204181"

$ws.Range("H17").Value = "1909"
$ws.Range("I17").Value = "1909"
$ws.Range("H18").Value = "1909"
$ws.Range("I18").Value = "1909"

$ws.Range("G17").Value = "posters"
$ws.Range("G18").Value = "posters"

$ws.Range("D17").Value = "Sports in art"
$ws.Range("D18").Value = "Sports in art"

$ws.Range("J1").Value = "publisher"
$ws.Range("J17").Value = "Norris & Winter"
$ws.Range("J18").Value = "Norris & Winter"

$ws.Range("E18").Value = "Baseball in art"
$ws.Range("E17").Value = "Football in art"

$ws.Range("C17").Value = "Poster of male football player"
$ws.Range("C18").Value = "Poster of male baseball player pitching"

$ws.Range("F1").Value = "subject3"
$ws.Range("F18").Value = "Tennessee Volunteers (Football team)"

$ws.Range("B17").Value = "Tennessee Volunteers Football poster"
$ws.Range("B18").Value = "Tennessee Volunteers Baseball poster"

$ws.Range("K1").Value = "geographic_subject"
$ws.Range("K17").Value = "Tennessee"
$ws.Range("K18").Value = "Tennessee"

# --- Column C (abstract) formatting: widen + wrap text --------------------
$ws.Columns.Item(3).ColumnWidth = 35.666666666666664
$ws.Range("C1:C18").WrapText = $true

# --- Column G (form, shifted) width ---------------------------------------
$ws.Columns.Item(7).ColumnWidth = 10.83

# --- Row heights for rows whose abstract text now wraps over several lines
$ws.Range("A3").RowHeight = 64
$ws.Range("A4").RowHeight = 64
$ws.Range("A5").RowHeight = 64
$ws.Range("A6").RowHeight = 64
$ws.Range("A9").RowHeight = 32
$ws.Range("A10").RowHeight = 80
$ws.Range("A12").RowHeight = 64
$ws.Range("A13").RowHeight = 64
$ws.Range("A14").RowHeight = 32
$ws.Range("A15").RowHeight = 32
$ws.Range("A16").RowHeight = 64

Write-Output "done"
